$d = $word.ActiveDocument

# Locate the run containing the closing ")}" field marker text.
$rng = $d.Content
$found = $rng.Find.Execute(")}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # $rng now covers ")}" exactly. Split it into two runs -- one for ")"
    # and one for "}" -- by toggling a character formatting property on
    # just the trailing "}" character. Word automatically breaks a run
    # whenever a sub-range ends up with different formatting than its
    # neighbour; flipping Bold on then off again leaves the visible
    # formatting unchanged but forces the run boundary to be created.
    $splitPoint = $rng.Start + 1
    $closingBrace = $d.Range($splitPoint, $rng.End)
    $closingBrace.Bold = 1
    $closingBrace.Bold = 0
}
